$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new rows (56-58) describing additional performance topics ---
# Values are entered in the same order the original author typed them
# (B56, C56, C57, B57) so that the shared-string table comes out identical.
$ws.Range("B56").Value = "partition view with UNIONs"
$ws.Range("C56").Value = "view with UNIONs"
$ws.Range("C57").Value = "nested views"
$ws.Range("B57").Value = "stored procedures"

# --- Formatting: copy the look of the existing "header" rows (B51:B55) onto B56:B57 ---
$ws.Range("B55").Copy()
$ws.Range("B56:B57").PasteSpecial(-4122)   # xlPasteFormats

# --- Formatting: copy the look of the existing value column (C45:C55) onto C56:C58 ---
$ws.Range("C55").Copy()
$ws.Range("C56:C58").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Match the taller row height used by the other header-style rows (47, 51-55)
$ws.Rows(56).RowHeight = 18
$ws.Rows(57).RowHeight = 18

# --- View state: scroll down and select the newly added cells ---
$ws.Range("C56:C58").Select()
$excel.ActiveWindow.ScrollRow = 33
